$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 26 de Marzo de 2020 a las 08:42'
$ws.Cells.Item(15, 2).Value = 5888
$ws.Cells.Item(15, 3).Value = 300
$ws.Cells.Item(15, 5).Value = 5845
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 8).Value = 34
$ws.Cells.Item(21, 4).Value = 6
$ws.Cells.Item(21, 5).Value = 2489
$ws.Cells.Item(25, 5).Value = 1576
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 21
$ws.Cells.Item(44, 2).Value = 678
$ws.Cells.Item(44, 3).Value = 21
$ws.Cells.Item(44, 5).Value = 622
$ws.Cells.Item(57, 4).Value = 190
$ws.Cells.Item(57, 5).Value = 225
$ws.Cells.Item(61, 5).Value = 363
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 6
$ws.Cells.Item(70, 1).Value = 'Taiwan'
$ws.Cells.Item(70, 2).Value = 252
$ws.Cells.Item(70, 3).Value = 17
$ws.Cells.Item(70, 4).Value = 29
$ws.Cells.Item(70, 5).Value = 221
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 8).Value = 2
$ws.Cells.Item(71, 1).Value = 'Bulgaria'
$ws.Cells.Item(71, 2).Value = 243
$ws.Cells.Item(71, 3).Value = 1
$ws.Cells.Item(71, 4).Value = 4
$ws.Cells.Item(71, 5).Value = 236
$ws.Cells.Item(71, 6).Value = 8
$ws.Cells.Item(71, 8).Value = 3
$ws.Cells.Item(95, 6).Value = 3
$ws.Cells.Item(101, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(101, 3).Value = 15
$ws.Cells.Item(101, 4).Value = 17
$ws.Cells.Item(101, 5).Value = 68
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 8).Value = 1
$ws.Cells.Item(102, 1).Value = 'Bielorrusia'
$ws.Cells.Item(102, 2).Value = 86
$ws.Cells.Item(102, 4).Value = 29
$ws.Cells.Item(102, 5).Value = 57
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(103, 1).Value = 'Afganistan'
$ws.Cells.Item(103, 2).Value = 84
$ws.Cells.Item(103, 4).Value = 2
$ws.Cells.Item(103, 5).Value = 80
$ws.Cells.Item(103, 8).Value = 2
$ws.Cells.Item(104, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(104, 2).Value = 80
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 3
$ws.Cells.Item(104, 5).Value = 77
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(105, 1).Value = 'Georgia'
$ws.Cells.Item(105, 2).Value = 77
$ws.Cells.Item(105, 3).Value = 2
$ws.Cells.Item(105, 4).Value = 10
$ws.Cells.Item(105, 5).Value = 67
$ws.Cells.Item(105, 6).Value = 1
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(106, 1).Value = 'Camerun'
$ws.Cells.Item(106, 2).Value = 75
$ws.Cells.Item(106, 4).Value = 2
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(107, 1).Value = 'Guadalupe'
$ws.Cells.Item(107, 2).Value = 73
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 72
$ws.Cells.Item(107, 6).Value = 4
$ws.Cells.Item(109, 1).Value = 'Montenegro'
$ws.Cells.Item(109, 2).Value = 67
$ws.Cells.Item(109, 3).Value = 14
$ws.Cells.Item(109, 5).Value = 66
$ws.Cells.Item(109, 6).Value = 1
$ws.Cells.Item(110, 1).Value = 'Martinica'
$ws.Cells.Item(110, 2).Value = 66
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 6).Value = 7
$ws.Cells.Item(110, 8).Value = 1
$ws.Cells.Item(111, 1).Value = 'Uzbekistan'
$ws.Cells.Item(111, 2).Value = 65
$ws.Cells.Item(111, 3).Value = 5
$ws.Cells.Item(111, 5).Value = 65
$ws.Cells.Item(111, 6).Value = 4
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(112, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(112, 2).Value = 60
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 59
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(113, 1).Value = 'Cuba'
$ws.Cells.Item(113, 2).Value = 57
$ws.Cells.Item(113, 4).Value = 1
$ws.Cells.Item(113, 5).Value = 55
$ws.Cells.Item(113, 6).Value = 2
$ws.Cells.Item(117, 1).Value = 'Mauricio'
$ws.Cells.Item(117, 6).Value = 1
$ws.Cells.Item(118, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(118, 6).Value = 0
